$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "61.678.30"
$ws.Range("E2").Value = "  -0.62%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.400.19"
$ws.Range("E3").Value = "  -0.42%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.22%  "

# Row 5 - BNB
$ws.Range("D5").Value = "411.85"
$ws.Range("E5").Value = "  +0.70%  "

# Row 6 - Solana
$ws.Range("D6").Value = "128.93"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.618"
$ws.Range("E7").Value = "  -2.79%  "

# Row 8 - USDC
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.720"
$ws.Range("E9").Value = "  -1.92%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -6.01%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "42.45"
$ws.Range("E11").Value = "  -0.06%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  +1.62%  "

# Row 13 - now WrappedliquidstakedEther2.0 (was ShibaInu)
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.939.15"
$ws.Range("E13").Value = "  -0.44%  "

# Row 14 - now ShibaInu (was WrappedliquidstakedEther2.0)
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000213"
$ws.Range("E14").Value = "  -3.46%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  -0.10%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "20.31"
$ws.Range("E16").Value = "  -2.29%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.392.10"
$ws.Range("E17").Value = "  -0.97%  "

# Row 18 - Uniswap
$ws.Range("D18").Value = "12.52"
$ws.Range("E18").Value = "  +3.20%  "

# Row 19 - Polygon
$ws.Range("E19").Value = "  -0.10%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "61.671.79"
$ws.Range("E20").Value = "  -0.44%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "479.29"
$ws.Range("E21").Value = "  +16.20%  "

# Row 22 - Litecoin
$ws.Range("D22").Value = "90.39"
$ws.Range("E22").Value = "  +1.42%  "

# Row 23 - ImmutableX
$ws.Range("D23").Value = "3.25"
$ws.Range("E23").Value = "  +2.52%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = "13.01"
$ws.Range("E24").Value = "  -0.52%  "

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "3.28"
$ws.Range("E25").Value = "  +1.43%  "

# Row 26 - Filecoin
$ws.Range("D26").Value = "9.69"
$ws.Range("E26").Value = "  +9.30%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "32.92"
$ws.Range("E27").Value = "  -0.85%  "

# Row 28 - LEO
$ws.Range("E28").Value = "  -0.30%  "

# Row 29 - RenderToken
$ws.Range("D29").Value = "7.66"
$ws.Range("E29").Value = "  +0.82%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -2.75%  "

# Row 31 - Cosmos
$ws.Range("D31").Value = "11.79"
$ws.Range("E31").Value = "  -1.10%  "

# Row 32 - Kaspa
$ws.Range("E32").Value = "  -2.14%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -3.64%  "

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = "40.74"
$ws.Range("E34").Value = "  -4.88%  "

# Row 35 - Dai
$ws.Range("E35").Value = "  -0.72%  "

# Row 36 - OKB
$ws.Range("D36").Value = "58.66"
$ws.Range("E36").Value = "  +8.18%  "

# Row 37 - VeChain
$ws.Range("D37").Value = "0.0483"
$ws.Range("E37").Value = "  -3.14%  "

# Row 38 - FirstDigitalUSD
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.24%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +3.59%  "

# Row 40 - Monero
$ws.Range("D40").Value = "149.04"
$ws.Range("E40").Value = "  +5.30%  "

# Row 41 - now Stellar (was TheGraph)
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.134"
$ws.Range("E41").Value = "  +0.16%  "

# Row 42 - now TheGraph (was Stellar)
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "0.320"
$ws.Range("E42").Value = "  +2.81%  "

# Row 43 - LidoDAOToken
$ws.Range("D43").Value = "3.32"
$ws.Range("E43").Value = "  -1.34%  "

# Row 44 - ARBITRUM
$ws.Range("D44").Value = "2.04"
$ws.Range("E44").Value = "  +3.82%  "

# Row 45 - WEMIXToken
$ws.Range("D45").Value = "2.59"
$ws.Range("E45").Value = "  +7.11%  "

# Row 46 - NEARProtocol
$ws.Range("D46").Value = "4.17"
$ws.Range("E46").Value = "  +1.64%  "

# Row 47 - ThetaToken
$ws.Range("E47").Value = "  +17.63%  "

# Row 48 - Celestia
$ws.Range("D48").Value = "16.31"
$ws.Range("E48").Value = "  -1.77%  "

# Row 49 - EnergySwap
$ws.Range("D49").Value = "22.08"
$ws.Range("E49").Value = "  +0.17%  "

# Row 50 - now BitcoinSV (was PEPE)
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "112.50"
$ws.Range("E50").Value = "  +14.16%  "

# Row 51 - now PEPE (was BitcoinSV)
$ws.Range("B51").Value = "PEPE"
$ws.Range("C51").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D51").Value = "0.0₃0510"
$ws.Range("E51").Value = "  +13.11%  "
